# Update task list: fill row 12 with the new "test question classification" task
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the task information for row 12 (previously empty)
$ws.Range("B12").Value = "Phân lọai các câu hỏi trong chương trình"
$ws.Range("C12").Value = "Hoang"
$ws.Range("D12").Value = "17/11"
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = "on processing"

# Row 12 needs extra height to fit the wrapped text
$ws.Rows.Item(12).RowHeight = 33

# Update the active selection to K12 (as recorded in the saved view)
$ws.Range("K12").Select()
